# Refresh the cryptos list (price / 1h-volume columns) as produced by the
# scheduled "Updated cryptos list ... with GitHub Actions" run, plus the
# PEPE / Bittensor row swap (rows 27-28) that came with this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These D-column price cells are plain numeric-looking strings (e.g. "1.00",
# "0.999", "549.75") that Excel would otherwise auto-convert to numbers on
# assignment. Force the cells to Text first so the values are stored
# verbatim, exactly like the source data (which keeps trailing zeros, etc.).
$textCells = @("D5","D6","D7","D8","D14","D18","D20","D25","D27","D29","D31","D32","D33","D36","D39","D41","D42","D45","D46","D47","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = '61.337.81'
$ws.Range("E2").Value = '  +0.18%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.380.45'
$ws.Range("E3").Value = '  +0.14%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.03%  '

# Row 5 - BNB
$ws.Range("D5").Value = '549.75'

# Row 6 - Solana
$ws.Range("D6").Value = '139.18'

# Row 7 - USDC
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.01%  '

# Row 8 - XRP
$ws.Range("D8").Value = '0.524'
$ws.Range("E8").Value = '  -0.99%  '

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = '2.380.52'
$ws.Range("E9").Value = '  +0.17%  '

# Row 10 - Dogecoin
$ws.Range("E10").Value = '  +2.13%  '

# Row 11 - TRON
$ws.Range("E11").Value = '  +1.44%  '

# Row 13 - Cardano
$ws.Range("E13").Value = '  +0.70%  '

# Row 14 - Avalanche
$ws.Range("D14").Value = '25.20'
$ws.Range("E14").Value = '  -0.43%  '

# Row 15 - ShibaInu
$ws.Range("E15").Value = '  +1.40%  '

# Row 16 - WrappedBTC
$ws.Range("D16").Value = '61.246.72'
$ws.Range("E16").Value = '  +0.18%  '

# Row 17 - WrappedEther
$ws.Range("D17").Value = '2.370.37'
$ws.Range("E17").Value = '  -0.32%  '

# Row 18 - Chainlink
$ws.Range("D18").Value = '10.97'
$ws.Range("E18").Value = '  +2.36%  '

# Row 19 - Polkadot
$ws.Range("E19").Value = '  +0.80%  '

# Row 20 - BitcoinCash
$ws.Range("D20").Value = '321.24'
$ws.Range("E20").Value = '  +0.96%  '

# Row 21 - Uniswap
$ws.Range("E21").Value = '  +0.84%  '

# Row 22 - Dai
$ws.Range("E22").Value = '  +0.06%  '

# Row 24 - SuiNetwork
$ws.Range("E24").Value = '  -8.97%  '

# Row 25 - Aptos
$ws.Range("D25").Value = '8.60'
$ws.Range("E25").Value = '  +5.03%  '

# Row 26 - InternetComputer(DFINITY)
$ws.Range("E26").Value = '  +1.48%  '

# Row 27 - was PEPE, now Bittensor
$ws.Range("B27").Value = 'Bittensor'
$ws.Range("C27").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D27").Value = '513.51'
$ws.Range("E27").Value = '  -2.85%  '

# Row 28 - was Bittensor, now PEPE
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0897'
$ws.Range("E28").Value = '  -3.08%  '

# Row 29 - Kaspa
$ws.Range("D29").Value = '0.151'
$ws.Range("E29").Value = '  +3.12%  '

# Row 30 - Fetch.AI
$ws.Range("E30").Value = '  -2.61%  '

# Row 31 - PancakeSwap
$ws.Range("D31").Value = '1.84'
$ws.Range("E31").Value = '  -0.05%  '

# Row 32 - ImmutableX
$ws.Range("D32").Value = '1.54'
$ws.Range("E32").Value = '  -2.21%  '

# Row 33 - FirstDigitalUSD
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.04%  '

# Row 34 - Stacks
$ws.Range("E34").Value = '  +3.57%  '

# Row 35 - NEARProtocol
$ws.Range("E35").Value = '  +0.90%  '

# Row 36 - RenderToken
$ws.Range("D36").Value = '5.47'
$ws.Range("E36").Value = '  +0.49%  '

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = '  +1.18%  '

# Row 38 - EthereumClassic
$ws.Range("E38").Value = '  +2.26%  '

# Row 39 - Monero
$ws.Range("D39").Value = '146.67'
$ws.Range("E39").Value = '  +4.32%  '

# Row 40 - USDe
$ws.Range("E40").Value = '  -0.12%  '

# Row 41 - OKB
$ws.Range("D41").Value = '41.27'
$ws.Range("E41").Value = '  +2.00%  '

# Row 42 - Aave
$ws.Range("D42").Value = '151.19'
$ws.Range("E42").Value = '  +7.58%  '

# Row 43 - dogwifhat
$ws.Range("E43").Value = '  +1.07%  '

# Row 44 - Filecoin
$ws.Range("E44").Value = '  -0.17%  '

# Row 45 - Hedera
$ws.Range("D45").Value = '0.0523'

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = '19.45'
$ws.Range("E46").Value = '  -3.23%  '

# Row 47 - Mantle
$ws.Range("D47").Value = '0.579'
$ws.Range("E47").Value = '  +0.86%  '

# Row 48 - Stellar
$ws.Range("D48").Value = '0.0907'
$ws.Range("E48").Value = '  +0.30%  '

# Row 49 - VeChain
$ws.Range("E49").Value = '  -0.29%  '

# Row 50 - WhiteBITCoin
$ws.Range("E50").Value = '  +0.44%  '

# Row 51 - EnergySwap
$ws.Range("D51").Value = '16.80'
$ws.Range("E51").Value = '  +0.22%  '
